# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -11
$ws.Range("F5").Value = -9
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -2
